# "vgg only svm run"
# Insert a new results row for an SVM run trained on VGG features.
# This new row lands at row 7 (pushing the existing XGBoost / Random
# Forest rows down by one), matching the target diff which turns the
# table from A1:F12 into A1:F13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 7.. down by inserting a fresh row at position 7.
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row with the SVM/VGG results.
$ws.Range("A7").Value = "SVM"
$ws.Range("B7").Value = "VGG"
$ws.Range("C7").Value = "C: 1, degree: 2, gamma: 0.001, kernel: rbf"
$ws.Range("D7").Value = 0.9431950691321006
$ws.Range("E7").Value = 0.8697316125451013
$ws.Range("F7").Value = 0.8643384822028207
